# Error Calculations and Plots
# Remove two rows (IDs "RM 232" and "SC 92") from the data table and update
# several previously-missing / now-missing cells in columns E ("D") and
# F ("F") to reflect the newly imputed / removed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the "RM 232" row (row 26) and the "SC 92" row (originally row
# 28, which becomes row 27 once the first row has been removed). Excel
# shifts everything below each deleted row up by one automatically. ---
$ws.Rows("26:26").Delete()
$ws.Rows("27:27").Delete()

# --- Update column F ("F") values that changed after the re-shuffle ---
$ws.Cells.Item(6, 6).Value = 16.43      # RM 21
$ws.Cells.Item(8, 6).ClearContents()    # RM 38
$ws.Cells.Item(18, 6).Value = 18.35     # RM 120
$ws.Cells.Item(20, 6).ClearContents()   # RM 134
$ws.Cells.Item(23, 6).Value = 16.48     # RM 140
$ws.Cells.Item(25, 6).ClearContents()   # RM 145
$ws.Cells.Item(30, 6).Value = 16.89     # SC 120

# --- Update column E ("D") values that changed after the re-shuffle ---
$ws.Cells.Item(27, 5).Value = -10       # SC 101
$ws.Cells.Item(28, 5).ClearContents()   # SC 105
$ws.Cells.Item(29, 5).ClearContents()   # SC 119
$ws.Cells.Item(30, 5).Value = -5.7      # SC 120
$ws.Cells.Item(32, 5).ClearContents()   # SC 193
